# Updating the staging env testdata
# Swap the "AddtionalParam" section identifiers between scenario1 (row 2)
# and scenario2 (row 6), and rotate the Expected_Study_Design_vals values
# for both scenario blocks (rows 2-4 and rows 6-8) accordingly. Also move
# the active selection to K6 and scroll the sheet so column E is leftmost.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- scenario1 block (rows 2-4): AddtionalParam section1 <-> section2 ---
$ws.Range("I2").Value2 = "study_design_section1"
$ws.Range("J2").Value2 = "study_design_section1_checkbox"

# --- scenario2 block (rows 6-8): AddtionalParam section2 <-> section1 ---
$ws.Range("I6").Value2 = "study_design_section2"
$ws.Range("J6").Value2 = "study_design_section2_checkbox"

# --- Expected_Study_Design_vals (column O) rotates within each block ---
# scenario1: O2 <- old O4, O3 <- old O2, O4 <- old O3
$ws.Range("O2").Value2 = "Phase 2 RCT"
$ws.Range("O3").Value2 = "Phase 3 RCT"
$ws.Range("O4").Value2 = "Phase NR RCT"

# scenario2: O6 <- old O8, O7 <- old O6, O8 <- old O7
$ws.Range("O6").Value2 = "Phase 2 RCT"
$ws.Range("O7").Value2 = "Phase 3 RCT"
$ws.Range("O8").Value2 = "Phase NR RCT"

# --- Update the view: move the active selection to K6, and scroll the
#     window so column E is the leftmost visible column (best-effort —
#     the window's own scroll-position attribute is session state that
#     this host does not persist independently of the frozen-pane split). ---
$ws.Range("K6").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
